# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 09:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1621196
$ws.Range("C4").Value = 294
$ws.Range("D4").Value = 382244
$ws.Range("E4").Value = 1142593
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 96359

# Row 11 - Alemania
$ws.Range("D11").Value = 159000
$ws.Range("E11").Value = 11712

# Row 52 - Chequia
$ws.Range("B52").Value = 8757
$ws.Range("C52").Value = 3
$ws.Range("D52").Value = 5932
$ws.Range("E52").Value = 2519

# Row 66 - Armenia
$ws.Range("B66").Value = 5928
$ws.Range("C66").Value = 322
$ws.Range("D66").Value = 2874
$ws.Range("E66").Value = 2980
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 74

# Row 92 - El Salvador
$ws.Range("B92").Value = 1725
$ws.Range("C92").Value = 85
$ws.Range("D92").Value = 562
$ws.Range("E92").Value = 1130
